$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column headers (row 1) stay the same.

# Data rows 2-7: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# followed by 16 numeric metric columns (E..T)
$rows = @(
    @{ Row = 2;  A = "ECs"; B = "Cdh1"; C = "Itgb7"; D = "ECs";
       Vals = @(1, 0.3333333333333333, 0.5650276666666666, 1.695083, 0.3735116729352592, 0.3735116729352592, 3, 1, 4.233575666666667, 12.700727, 0.4702904532725788, 0.4702904532725788, 2.392087380593444, 21.528786425341, 0.1756589739673222, 0.1756589739673222) },
    @{ Row = 3;  A = "ECs"; B = "Cdh1"; C = "Itgb7"; D = "FAPs";
       Vals = @(1, 0.3333333333333333, 0.5650276666666666, 1.695083, 0.3735116729352592, 0.3735116729352592, 3, 1, 3.605649999999999, 10.81695, 0.4005367817548413, 0.4005367817548413, 2.037292006316666, 18.33562805684999, 0.1496051634253555, 0.1496051634253555) },
    @{ Row = 4;  A = "ECs"; B = "Cdh1"; C = "Itgb7"; D = "sCs";
       Vals = @(1, 0.3333333333333333, 0.5650276666666666, 1.695083, 0.3735116729352592, 0.3735116729352592, 3, 1, 1.162819, 3.488456999999999, 0.1291727649725799, 0.1291727649725799, 0.6570249063256665, 5.913224156930998, 0.04824753554258138, 0.04824753554258138) },
    @{ Row = 5;  A = "sCs"; B = "Cdh1"; C = "Itgb7"; D = "ECs";
       Vals = @(3, 1, 0.9477166666666665, 2.84315, 0.6264883270647409, 0.6264883270647408, 3, 1, 4.233575666666667, 12.700727, 0.4702904532725788, 0.4702904532725788, 4.012230218894444, 36.11007197004999, 0.2946314793052566, 0.2946314793052565) },
    @{ Row = 6;  A = "sCs"; B = "Cdh1"; C = "Itgb7"; D = "FAPs";
       Vals = @(3, 1, 0.9477166666666665, 2.84315, 0.6264883270647409, 0.6264883270647408, 3, 1, 3.605649999999999, 10.81695, 0.4005367817548413, 0.4005367817548413, 3.417134599166666, 30.75421139249999, 0.2509316183294857, 0.2509316183294857) },
    @{ Row = 7;  A = "sCs"; B = "Cdh1"; C = "Itgb7"; D = "sCs";
       Vals = @(3, 1, 0.9477166666666665, 2.84315, 0.6264883270647409, 0.6264883270647408, 3, 1, 1.162819, 3.488456999999999, 0.1291727649725799, 0.1291727649725799, 1.102022946616666, 9.918206519549997, 0.08092522942999857, 0.08092522942999855) }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D

    $col = 5
    foreach ($v in $r.Vals) {
        $ws.Cells.Item($rowNum, $col).Value = $v
        $col = $col + 1
    }
}
